$p = $ppt.ActivePresentation

# --- Slide 1: fix title wording ("reduce death" -> "reduce death rates") ---
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(5)
$titleTr = $titleShape.TextFrame.TextRange
$titlePara = $titleTr.Paragraphs(1, 1)
$titleRun = $titlePara.Runs(1, 1)
$titleRun.Text = "Project 1: What can we learn about COVID-19 trends to reduce death rates in a future pandemic?"

# --- Slide 2: Data and Aims bullet updates ---
$s2 = $p.Slides.Item(2)
$bodyShape = $s2.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange

# Paragraph 2: worldometers sentence gets an extra clause appended
$para2 = $bodyTr.Paragraphs(2, 1)
$run2 = $para2.Runs(1, 1)
$run2.Text = "The CSV dataset was pulled from the worldometers.info website and made available at kaggle.com"

# Paragraph 3: country count 255 -> 225
$para3 = $bodyTr.Paragraphs(3, 1)
$run3 = $para3.Runs(1, 1)
$run3.Text = "Included is daily data for 225 countries on case counts and death counts"

# Paragraph 6: empty paragraph right after "Aims" gets removed
$para6 = $bodyTr.Paragraphs(6, 1)
$para6.Delete()
